$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPV")
$tbl = $ws.ListObjects.Item("ppv")

# Add two new columns to the "ppv" table: Program, Lot
$col9 = $tbl.ListColumns.Add()
$ws.Range("I1").Value = "Program"

$col10 = $tbl.ListColumns.Add()
$ws.Range("J1").Value = "Lot"

$ws.Range("J2").Select()
